$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- MatrNr value: was a placeholder "?" text, now the actual number ---
$ws.Range("B2").Value = 6003356

# --- Team value: was a placeholder "?" text, now the team/project name ---
$ws.Range("B4").Value = "RayForge"

# --- Add a new time-tracking entry (row 13) for the meeting on 2024-10-28 ---

# Match formatting of the existing rows (copy formats only, so the shared
# style indices used by the sheet are reused instead of new ones minted).
$ws.Range("A7").Copy() | Out-Null
$ws.Range("A13").PasteSpecial(-4122) | Out-Null

$ws.Range("B7").Copy() | Out-Null
$ws.Range("B13").PasteSpecial(-4122) | Out-Null

$ws.Range("C9").Copy() | Out-Null
$ws.Range("C13").PasteSpecial(-4122) | Out-Null

$ws.Range("D9").Copy() | Out-Null
$ws.Range("D13").PasteSpecial(-4122) | Out-Null

$ws.Application.CutCopyMode = 0

$ws.Range("A13").Value = 45593
$ws.Range("B13").Value = 2
$ws.Range("C13").Value = "Besprechung"
$ws.Range("D13").Value = "Aufgabenverteilung, Organisierung "

# --- Update the active selection to the newly filled cell ---
$ws.Range("D13").Select() | Out-Null
